# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Columns D (Price) and E (Volume(1h)) hold human-formatted text, not real
# numbers (e.g. "65.169.65", "  +0.37%  "), so every value is written with a
# leading apostrophe to force text entry the same way a user typing into
# Excel would - this avoids Excel reinterpreting them as numbers/dates and
# silently mangling values such as "1.00" -> 1 or "0.0000183" -> 1.83E-05.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.169.65"
$ws.Range("E2").Value = "'  +0.37%  "
$ws.Range("D3").Value = "'3.540.90"
$ws.Range("E3").Value = "'  +4.45%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'598.77"
$ws.Range("E5").Value = "'  +3.66%  "
$ws.Range("D6").Value = "'138.29"
$ws.Range("E6").Value = "'  +2.42%  "
$ws.Range("D7").Value = "'3.540.90"
$ws.Range("E7").Value = "'  +4.51%  "
$ws.Range("E8").Value = "'  +0.10%  "
$ws.Range("E9").Value = "'  +0.68%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("D11").Value = "'6.92"
$ws.Range("E11").Value = "'  -2.02%  "
$ws.Range("E12").Value = "'  +4.49%  "
$ws.Range("D13").Value = "'4.142.24"
$ws.Range("E13").Value = "'  +4.44%  "
$ws.Range("D14").Value = "'0.0000183"
$ws.Range("E14").Value = "'  +3.87%  "
$ws.Range("D15").Value = "'27.37"
$ws.Range("E15").Value = "'  +5.40%  "
$ws.Range("D16").Value = "'3.543.36"
$ws.Range("E16").Value = "'  +4.35%  "
$ws.Range("D18").Value = "'65.106.81"
$ws.Range("E18").Value = "'  +0.32%  "
$ws.Range("E19").Value = "'  +6.43%  "
$ws.Range("D20").Value = "'5.87"
$ws.Range("E20").Value = "'  +1.38%  "
$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "'  +5.87%  "
$ws.Range("D22").Value = "'391.73"
$ws.Range("E22").Value = "'  +3.45%  "
$ws.Range("D23").Value = "'0.574"
$ws.Range("E23").Value = "'  +4.85%  "
$ws.Range("D24").Value = "'3.682.26"
$ws.Range("E24").Value = "'  +4.44%  "
$ws.Range("D25").Value = "'73.85"
$ws.Range("E25").Value = "'  +2.97%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "'  -0.09%  "
$ws.Range("E27").Value = "'  +10.21%  "
$ws.Range("D28").Value = "'7.85"
$ws.Range("E28").Value = "'  +12.98%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "'  +4.19%  "
$ws.Range("D31").Value = "'8.33"
$ws.Range("E31").Value = "'  +5.02%  "
$ws.Range("D32").Value = "'3.563.24"
$ws.Range("E32").Value = "'  +4.81%  "
$ws.Range("D33").Value = "'1.38"
$ws.Range("E33").Value = "'  +22.20%  "
$ws.Range("E34").Value = "'  +0.03%  "
$ws.Range("D35").Value = "'23.84"
$ws.Range("E35").Value = "'  +5.20%  "
$ws.Range("E36").Value = "'  +2.96%  "
$ws.Range("D37").Value = "'1.59"
$ws.Range("E37").Value = "'  +9.71%  "

# Monero/Aptos swapped ranking positions (row 38 <-> row 39); the index
# column (A) stays put, only Coin/Link/Price/Volume move.
$ws.Range("B38").Value = "'Monero"
$ws.Range("C38").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'170.33"
$ws.Range("E38").Value = "'  +1.14%  "
$ws.Range("B39").Value = "'Aptos"
$ws.Range("C39").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'6.94"
$ws.Range("E39").Value = "'  +4.99%  "
$ws.Range("D40").Value = "'5.03"
$ws.Range("E40").Value = "'  +8.80%  "
$ws.Range("D41").Value = "'0.0806"
$ws.Range("E41").Value = "'  +7.86%  "
$ws.Range("D42").Value = "'0.824"
$ws.Range("E42").Value = "'  +2.15%  "
$ws.Range("D43").Value = "'26.55"
$ws.Range("E43").Value = "'  +22.59%  "
$ws.Range("D44").Value = "'42.55"
$ws.Range("E44").Value = "'  -1.78%  "
$ws.Range("E45").Value = "'  -0.13%  "
$ws.Range("E46").Value = "'  +3.64%  "
$ws.Range("E47").Value = "'  +10.93%  "
$ws.Range("E48").Value = "'  +5.48%  "
$ws.Range("D49").Value = "'6.84"
$ws.Range("E49").Value = "'  +6.64%  "
$ws.Range("D50").Value = "'2.411.90"
$ws.Range("E50").Value = "'  +12.31%  "
$ws.Range("D51").Value = "'311.11"
$ws.Range("E51").Value = "'  +18.52%  "
